$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 9: I1_R1 / Reale (penalità 5)
$ws.Range("A9").Value = "I1_R1"
# New row 10: I1_R2 / Reale (penalità 20) + capacità 1000 + service time 120
$ws.Range("A10").Value = "I1_R2"
# Write column B after column A so the shared-string table is built in the
# same order id-then-description for each row (matches target sharedStrings order)
$ws.Range("B9").Value = "Reale (penalità 5)"
$ws.Range("B10").Value = "Reale (penalità 20) + capacità 1000 + service time 120"

# Row 11: an empty, bolded cell - user prepped the next row with Bold
# formatting (this registers a new cell style / cellXfs entry).
$ws.Range("B11").Font.Bold = $true

# Move the active selection to the new empty row, mirroring where the
# user's cursor ended up after entering the new data.
$ws.Range("B11").Select() | Out-Null
